$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing table (Tabla1) to include the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("D1:I30"))

# Fill in the new test-plan rows (16-30) covering "Crear Modelo" functional tests
$ws.Range("D16").Value = 'CrM13'
$ws.Range("E16").Value = 'Preprocesamiento de Datos'
$ws.Range("F16").Value = 'CrM11||CrM10'
$ws.Range("G16").Value = '1-Se selecciona Eliminar filas y se aplica'
$ws.Range("H16").Value = 'Barra de carga y Notificación de Éxito'
$ws.Range("I16").Value = '✅'
$ws.Range("D17").Value = 'CrM14'
$ws.Range("E17").Value = 'Preprocesamiento de Datos'
$ws.Range("F17").Value = 'CrM11||CrM10'
$ws.Range("G17").Value = '1-Se selecciona Rellenar con Media y se Aplica'
$ws.Range("H17").Value = 'Barra de carga y Notificación de Éxito'
$ws.Range("I17").Value = '✅'
$ws.Range("D18").Value = 'CrM15'
$ws.Range("E18").Value = 'Preprocesamiento de Datos'
$ws.Range("F18").Value = 'CrM11||CrM10'
$ws.Range("G18").Value = '1-Se selecciona Rellenar con Mediana y se Aplica'
$ws.Range("H18").Value = 'Barra de carga y Notificación de Éxito'
$ws.Range("I18").Value = '✅'
$ws.Range("D19").Value = 'CrM16'
$ws.Range("E19").Value = 'Preprocesamiento de Datos'
$ws.Range("F19").Value = 'CrM11||CrM10'
$ws.Range("G19").Value = '1-Se selecciona Rellenar con Constante y se Aplica'
$ws.Range("H19").Value = 'Notificación Error de Validación'
$ws.Range("I19").Value = '✅'
$ws.Range("D20").Value = 'CrM17'
$ws.Range("E20").Value = 'Preprocesamiento de Datos'
$ws.Range("F20").Value = 'CrM11||CrM10'
$ws.Range("G20").Value = '1-Se selecciona Rellenar con Constante, se escribe una constante y Aplica'
$ws.Range("H20").Value = 'Barra de carga y Notificación de Éxito'
$ws.Range("I20").Value = '✅'
$ws.Range("D21").Value = 'CrM18'
$ws.Range("E21").Value = 'Separación de Datos'
$ws.Range("F21").Value = 'CrM15||13||14||17||8'
$ws.Range("G21").Value = '1-Se acepta o cierra la notficación'
$ws.Range("H21").Value = 'Nueva sección División de Datos'
$ws.Range("I21").Value = '✅'
$ws.Range("D22").Value = 'CrM19'
$ws.Range("E22").Value = 'Separación de Datos'
$ws.Range("F22").Value = 'CrM18'
$ws.Range("G22").Value = '1-Se escoge un valor no numérico de Semilla y se Divide'
$ws.Range("H22").Value = 'Notificación de Error'
$ws.Range("I22").Value = '✅'
$ws.Range("D23").Value = 'CrM20'
$ws.Range("E23").Value = 'Separación de Datos'
$ws.Range("F23").Value = 'CrM18'
$ws.Range("G23").Value = '1- No se escoge un valor entre 0 y 10000000000 y se Divide'
$ws.Range("H23").Value = 'Notificación de Error'
$ws.Range("I23").Value = '❌'
$ws.Range("D24").Value = 'CrM21'
$ws.Range("E24").Value = 'Separación de Datos'
$ws.Range("F24").Value = 'CrM18'
$ws.Range("G24").Value = '1-Se escoge un valor entre 0 y 10000000000 y se Divide'
$ws.Range("H24").Value = 'Notificación de Éxito'
$ws.Range("I24").Value = '❌'
$ws.Range("D25").Value = 'CrM22'
$ws.Range("E25").Value = 'Separación de Datos'
$ws.Range("F25").Value = 'CrM18'
$ws.Range("G25").Value = '1-Se comprueban la funcionalidad del slider y se Divide'
$ws.Range("H25").Value = 'Notificación de Éxito'
$ws.Range("I25").Value = '❌'
$ws.Range("D26").Value = 'CrM23'
$ws.Range("E26").Value = 'Crear Modelo'
$ws.Range("F26").Value = 'CrM21||CrM22'
$ws.Range("G26").Value = '1-Se acepta o cierra la notficación'
$ws.Range("H26").Value = 'Nueva sección Visualización dle Modelo'
$ws.Range("I26").Value = '✅'
$ws.Range("D27").Value = 'CrM24'
$ws.Range("E27").Value = 'Crear Modelo'
$ws.Range("F27").Value = 'CrM23 && CrM8 '
$ws.Range("G27").Value = '1-Se acepta la notificación de más de una variable de entrada'
$ws.Range("H27").Value = 'Nueva sección Visualización dle Modelo'
$ws.Range("I27").Value = '✅'
$ws.Range("D28").Value = 'CrM25'
$ws.Range("E28").Value = 'Crear Modelo'
$ws.Range("F28").Value = 'CrM23 ||CrM24'
$ws.Range("G28").Value = '1-Se comprueba la fórmula con los nombres de las columnas seleccionadas'
$ws.Range("H28").Value = 'Presentación correcta de la fórmula'
$ws.Range("I28").Value = '✅'
$ws.Range("D29").Value = 'CrM26'
$ws.Range("E29").Value = 'Crear Modelo'
$ws.Range("F29").Value = 'CrM23||CrM24'
$ws.Range("G29").Value = '1-Se comprueba R 2 y el EMC de entrenamiento y test'
$ws.Range("H29").Value = 'Valores correctamente calculados'
$ws.Range("I29").Value = '✅'
$ws.Range("D30").Value = 'CrM27'
$ws.Range("E30").Value = 'Crear Modelo'
$ws.Range("F30").Value = '(CrM23||CrM24)&&CrM10'
$ws.Range("G30").Value = '1-Se comprueban ambas gráficas'
$ws.Range("H30").Value = 'Representación correccta '
$ws.Range("I30").Value = '✅'

# Leave selection on F16, matching the final cursor position after data entry
$ws.Range("F16").Select()
